$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.08743719878012257
$ws.Range("D2").Value = 0.6279230237380489
$ws.Range("E2").Value = 0.06789468549853694
$ws.Range("F2").Value = 17.98820901680062
$ws.Range("G2").Value = 0.003146068211839157
$ws.Range("J2").Value = 0.2197593822912154
$ws.Range("L2").Value = 0.2335861272159718
$ws.Range("M2").Value = 16.22287908741822

$ws.Range("C3").Value = 0.08899340847368364
$ws.Range("D3").Value = 0.6283865045828634
$ws.Range("E3").Value = 0.05882268558075054
$ws.Range("F3").Value = 18.4134106178484
$ws.Range("G3").Value = 0.003171861845197182
$ws.Range("J3").Value = 0.2260456287015842
$ws.Range("L3").Value = 0.2085333253586725
$ws.Range("M3").Value = 15.25129885501514

$ws.Range("C4").Value = 0.09000477119064243
$ws.Range("D4").Value = 0.6296260947131316
$ws.Range("E4").Value = 0.05327544463061429
$ws.Range("F4").Value = 18.69681349105085
$ws.Range("G4").Value = 0.003188345837087719
$ws.Range("J4").Value = 0.2302370251315224
$ws.Range("L4").Value = 0.1932644321379371
$ws.Range("M4").Value = 14.65993981377983

$ws.Range("C5").Value = 0.09043094390248285
$ws.Range("D5").Value = 0.6303649232893633
$ws.Range("E5").Value = 0.05101969978699827
$ws.Range("F5").Value = 18.81783857041657
$ws.Range("G5").Value = 0.003195227791812577
$ws.Range("J5").Value = 0.2320272487579231
$ws.Range("L5").Value = 0.1870678378296589
$ws.Range("M5").Value = 14.42015991544895

$ws.Range("C6").Value = 0.0905025570073974
$ws.Range("D6").Value = 0.6305015359305344
$ws.Range("E6").Value = 0.05064539746499008
$ws.Range("F6").Value = 18.83826683632213
$ws.Range("G6").Value = 0.003196380533390744
$ws.Range("J6").Value = 0.2323294451940185
$ws.Range("L6").Value = 0.1860403611075441
$ws.Range("M6").Value = 14.38041448382387

$ws.Range("C7").Value = 0.09001046188993911
$ws.Range("D7").Value = 0.6296351209490467
$ws.Range("E7").Value = 0.05324500468328353
$ws.Range("F7").Value = 18.69842336288465
$ws.Range("G7").Value = 0.003188437980751565
$ws.Range("J7").Value = 0.2302608374177382
$ws.Range("L7").Value = 0.19318076282471
$ws.Range("M7").Value = 14.6567013148786

$ws.Range("C8").Value = 0.08796218379775667
$ws.Range("D8").Value = 0.6278806367612333
$ws.Range("E8").Value = 0.06476127157414879
$ws.Range("F8").Value = 18.13013630857839
$ws.Range("G8").Value = 0.003154828886542439
$ws.Range("J8").Value = 0.2218573096548582
$ws.Range("L8").Value = 0.2249224867934743
$ws.Range("M8").Value = 15.88673564556188

$ws.Range("C9").Value = 0.0843891023758232
$ws.Range("D9").Value = 0.6323251919796462
$ws.Range("E9").Value = 0.08757576848373816
$ws.Range("F9").Value = 17.1965825070094
$ws.Range("G9").Value = 0.003093956939485412
$ws.Range("J9").Value = 0.2080673511151687
$ws.Range("L9").Value = 0.2882156917437442
$ws.Range("M9").Value = 18.34517457801371

$ws.Range("C10").Value = 0.08203527411139788
$ws.Range("D10").Value = 0.640850854790898
$ws.Range("E10").Value = 0.1045525567301056
$ws.Range("F10").Value = 16.62636909177962
$ws.Range("G10").Value = 0.003052164581742289
$ws.Range("J10").Value = 0.1996631104378253
$ws.Range("L10").Value = 0.3355780178573866
$ws.Range("M10").Value = 20.18752276237825

$ws.Range("C11").Value = 0.08102365628518271
$ws.Range("D11").Value = 0.6459792490450127
$ws.Range("E11").Value = 0.1123403550260775
$ws.Range("F11").Value = 16.39336668232761
$ws.Range("G11").Value = 0.003033755955177947
$ws.Range("J11").Value = 0.1962364461649884
$ws.Range("L11").Value = 0.3573649268449799
$ws.Range("M11").Value = 21.03542539613323

$ws.Range("C12").Value = 0.08064912265164992
$ws.Range("D12").Value = 0.6481106248887443
$ws.Range("E12").Value = 0.1153003244797901
$ws.Range("F12").Value = 16.30904609201025
$ws.Range("G12").Value = 0.00302686897900129
$ws.Range("J12").Value = 0.1949978697167154
$ws.Range("L12").Value = 0.3656545512791638
$ws.Range("M12").Value = 21.35808603864331

$ws.Range("C13").Value = 0.08072940449446264
$ws.Range("D13").Value = 0.6476430125666184
$ws.Range("E13").Value = 0.1146623312440269
$ws.Range("F13").Value = 16.32703007614009
$ws.Range("G13").Value = 0.00302834852269064
$ws.Range("J13").Value = 0.1952619600784402
$ws.Range("L13").Value = 0.3638674021230486
$ws.Range("M13").Value = 21.28852227836404

$ws.Range("C14").Value = 0.08099267170075564
$ws.Range("D14").Value = 0.6461507459871711
$ws.Range("E14").Value = 0.1125836457024221
$ws.Range("F14").Value = 16.38635046431341
$ws.Range("G14").Value = 0.003033187692857474
$ws.Range("J14").Value = 0.1961333525687081
$ws.Range("L14").Value = 0.358046101711011
$ws.Range("M14").Value = 21.06193823862844

$ws.Range("C15").Value = 0.08115504439716403
$ws.Range("D15").Value = 0.6452616473193586
$ws.Range("E15").Value = 0.1113118584831483
$ws.Range("F15").Value = 16.42319904910079
$ws.Range("G15").Value = 0.003036162675726745
$ws.Range("J15").Value = 0.1966748555264388
$ws.Range("L15").Value = 0.3544856620093242
$ws.Range("M15").Value = 20.92335973314061

$ws.Range("C16").Value = 0.08210257911244412
$ws.Range("D16").Value = 0.6405416622762061
$ws.Range("E16").Value = 0.1040450434703288
$ws.Range("F16").Value = 16.64213724310696
$ws.Range("G16").Value = 0.003053379531574428
$ws.Range("J16").Value = 0.1998951972402239
$ws.Range("L16").Value = 0.334159442900841
$ws.Range("M16").Value = 20.1323217435413

$ws.Range("C17").Value = 0.08269903621307151
$ws.Range("D17").Value = 0.6379728958787609
$ws.Range("E17").Value = 0.09960479805659617
$ws.Range("F17").Value = 16.78328896532111
$ws.Range("G17").Value = 0.003064094047709098
$ws.Range("J17").Value = 0.2019736966631598
$ws.Range("L17").Value = 0.3217550359560732
$ws.Range("M17").Value = 19.64967297315059

$ws.Range("C18").Value = 0.08304766744247871
$ws.Range("D18").Value = 0.6366124539842701
$ws.Range("E18").Value = 0.09705690156864932
$ws.Range("F18").Value = 16.86695220185675
$ws.Range("G18").Value = 0.003070313676818434
$ws.Range("J18").Value = 0.2032063811138514
$ws.Range("L18").Value = 0.3146427850496991
$ws.Range("M18").Value = 19.37297860478697

$ws.Range("C19").Value = 0.08316666268860473
$ws.Range("D19").Value = 0.6361716522414724
$ws.Range("E19").Value = 0.09619521527004338
$ws.Range("F19").Value = 16.89570103255807
$ws.Range("G19").Value = 0.003072429393188574
$ws.Range("J19").Value = 0.2036300751360685
$ws.Range("L19").Value = 0.312238416304325
$ws.Range("M19").Value = 19.2794464705334

$ws.Range("C20").Value = 0.0826349661357284
$ws.Range("D20").Value = 0.6382341642896279
$ws.Range("E20").Value = 0.1000768363207172
$ws.Range("F20").Value = 16.76800599093048
$ws.Range("G20").Value = 0.003062947597123266
$ws.Range("J20").Value = 0.2017485741422504
$ws.Range("L20").Value = 0.3230731516077583
$ws.Range("M20").Value = 19.70095615351494

$ws.Range("C21").Value = 0.08091511152115061
$ws.Range("D21").Value = 0.6465838415677752
$ws.Range("E21").Value = 0.1131938972142734
$ws.Range("F21").Value = 16.36881947015218
$ws.Range("G21").Value = 0.003031764054960654
$ws.Range("J21").Value = 0.1958757845521006
$ws.Range("L21").Value = 0.3597548506344026
$ws.Range("M21").Value = 21.12844722178318

$ws.Range("C22").Value = 0.07984091077118194
$ws.Range("D22").Value = 0.6531495609804381
$ws.Range("E22").Value = 0.1218309966167936
$ws.Range("F22").Value = 16.13079295838844
$ws.Range("G22").Value = 0.003011871868923239
$ws.Range("J22").Value = 0.1923826970093359
$ws.Range("L22").Value = 0.3839605144686686
$ws.Range("M22").Value = 22.07068290777755

$ws.Range("C23").Value = 0.08040965879411743
$ws.Range("D23").Value = 0.6495405077288297
$ws.Range("E23").Value = 0.117214779942536
$ws.Range("F23").Value = 16.25569900531696
$ws.Range("G23").Value = 0.003022444990320329
$ws.Range("J23").Value = 0.1942147320466034
$ws.Range("L23").Value = 0.3710186269556459
$ws.Range("M23").Value = 21.5668850199857

$ws.Range("C24").Value = 0.0826639144049679
$ws.Range("D24").Value = 0.6381156830744033
$ws.Range("E24").Value = 0.09986341280060174
$ws.Range("F24").Value = 16.77490760428003
$ws.Range("G24").Value = 0.003063465721061782
$ws.Range("J24").Value = 0.2018502346712694
$ws.Range("L24").Value = 0.3224771720411184
$ws.Range("M24").Value = 19.67776858935935

$ws.Range("C25").Value = 0.08530815092014876
$ws.Range("D25").Value = 0.6302335768824037
$ws.Range("E25").Value = 0.08137178122393607
$ws.Range("F25").Value = 17.42923342300048
$ws.Range("G25").Value = 0.003109899322911495
$ws.Range("J25").Value = 0.2115013928380947
$ws.Range("L25").Value = 0.2709590750584425
$ws.Range("M25").Value = 17.67436448148129
